# Add 4 new marcaciones rows (26-29) to the Marcaciones sheet, matching
# the data appended in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rows: nombre, cedula, agencia, horaEntrada, horaSalida, observaciones
$newRows = @(
    @("jamilton",   "7878",       "Alameda",  "22/05/2025, 14:57:09", "",                       "hou"),
    @("David",      "1265",       "Neiva",    "29/05/2025, 09:50:42", "",                       "aaa"),
    @("Goliat",     "14569",      "Cafetero", "29/05/2025, 10:06:14", "",                       "Llegue tarde"),
    @("Juan Pablo",  "1006036679", "Delicias", "28/05/2025, 10:07:46", "28/05/2025, 10:07:46",   "HGola")
)

$startRow = 26

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    for ($c = 1; $c -le 6; $c++) {
        $value = $rowData[$c - 1]

        # Skip blank "horaSalida" cells entirely - leave them empty, same
        # as the rest of the sheet's empty-string cells.
        if ($value -eq "") {
            continue
        }

        $cell = $ws.Cells.Item($r, $c)

        # Force text storage for numeric-looking cedula values so they are
        # kept as text (matching t="str" in the original file) instead of
        # being coerced into numbers.
        $numericLooking = $value -match '^[0-9]+$'
        if ($numericLooking) {
            $cell.NumberFormat = "@"
        }

        $cell.Value = $value
    }
}
